$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows append below the existing data (rows 4-7).
# Quantity / Cost Per / Total Cost columns mirror the existing sheet's
# convention of storing numeric-looking values as text, so force the
# number format to Text before writing them.
$ws.Range("C4:E7").NumberFormat = "@"

$ws.Range("A4").Value = "wholeg"
$ws.Range("B4").Value = "Milk - Whole"
$ws.Range("C4").Value = "9"
$ws.Range("D4").Value = "16.08"
$ws.Range("E4").Value = "144.72"

$ws.Range("A5").Value = "skimg"
$ws.Range("B5").Value = "Milk - Skim"
$ws.Range("C5").Value = "2"
$ws.Range("D5").Value = "15.56"
$ws.Range("E5").Value = "31.12"

$ws.Range("A6").Value = "twog"
$ws.Range("B6").Value = "Milk - 2%"
$ws.Range("C6").Value = "20"
$ws.Range("D6").Value = "16.08"
$ws.Range("E6").Value = "321.60"

$ws.Range("A7").Value = "choqt"
$ws.Range("B7").Value = "Milk - Chocolate (9/32oz)"
$ws.Range("C7").Value = "3"
$ws.Range("D7").Value = "16.65"
$ws.Range("E7").Value = "49.95"
